$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D3" = "53af4926-52ee-41d0-9acc-ae7230400002"
    "D4" = "53af4926-52ee-41d0-9acc-ae7230400003"
    "D5" = "53af4926-52ee-41d0-9acc-ae7230400004"
    "D6" = "53af4926-52ee-41d0-9acc-ae7230400005"
    "D7" = "53af4926-52ee-41d0-9acc-ae7230400001"
    "D9" = "53af4926-52ee-41d0-9acc-ae7230400003"
    "D10" = "53af4926-52ee-41d0-9acc-ae7230400004"
    "D11" = "53af4926-52ee-41d0-9acc-ae7230400005"
    "D12" = "53af4926-52ee-41d0-9acc-ae7230400001"
    "D13" = "53af4926-52ee-41d0-9acc-ae7230400002"
    "D14" = "53af4926-52ee-41d0-9acc-ae7230400003"
    "D16" = "53af4926-52ee-41d0-9acc-ae7230400005"
    "D17" = "53af4926-52ee-41d0-9acc-ae7230400001"
    "D18" = "53af4926-52ee-41d0-9acc-ae7230400002"
    "D19" = "53af4926-52ee-41d0-9acc-ae7230400003"
    "D20" = "53af4926-52ee-41d0-9acc-ae7230400004"
    "D23" = "53af4926-52ee-41d0-9acc-ae7230400007"
    "D24" = "53af4926-52ee-41d0-9acc-ae7230400008"
    "D25" = "53af4926-52ee-41d0-9acc-ae7230400009"
    "D26" = "53af4926-52ee-41d0-9acc-ae7230400010"
    "D27" = "53af4926-52ee-41d0-9acc-ae7230400006"
    "D29" = "53af4926-52ee-41d0-9acc-ae7230400008"
    "D30" = "53af4926-52ee-41d0-9acc-ae7230400009"
    "D31" = "53af4926-52ee-41d0-9acc-ae7230400010"
    "D32" = "53af4926-52ee-41d0-9acc-ae7230400006"
    "D33" = "53af4926-52ee-41d0-9acc-ae7230400007"
    "D34" = "53af4926-52ee-41d0-9acc-ae7230400008"
    "D36" = "53af4926-52ee-41d0-9acc-ae7230400010"
    "D37" = "53af4926-52ee-41d0-9acc-ae7230400006"
    "D38" = "53af4926-52ee-41d0-9acc-ae7230400007"
    "D39" = "53af4926-52ee-41d0-9acc-ae7230400008"
    "D40" = "53af4926-52ee-41d0-9acc-ae7230400009"
    "D43" = "53af4926-52ee-41d0-9acc-ae7230400012"
    "D44" = "53af4926-52ee-41d0-9acc-ae7230400013"
    "D45" = "53af4926-52ee-41d0-9acc-ae7230400014"
    "D46" = "53af4926-52ee-41d0-9acc-ae7230400015"
    "D47" = "53af4926-52ee-41d0-9acc-ae7230400011"
    "D49" = "53af4926-52ee-41d0-9acc-ae7230400013"
    "D50" = "53af4926-52ee-41d0-9acc-ae7230400014"
    "D51" = "53af4926-52ee-41d0-9acc-ae7230400015"
    "D52" = "53af4926-52ee-41d0-9acc-ae7230400011"
    "D53" = "53af4926-52ee-41d0-9acc-ae7230400012"
    "D54" = "53af4926-52ee-41d0-9acc-ae7230400013"
    "D56" = "53af4926-52ee-41d0-9acc-ae7230400015"
    "D57" = "53af4926-52ee-41d0-9acc-ae7230400011"
    "D58" = "53af4926-52ee-41d0-9acc-ae7230400012"
    "D59" = "53af4926-52ee-41d0-9acc-ae7230400013"
    "D60" = "53af4926-52ee-41d0-9acc-ae7230400014"
    "D63" = "53af4926-52ee-41d0-9acc-ae7230400017"
    "D64" = "53af4926-52ee-41d0-9acc-ae7230400018"
    "D65" = "53af4926-52ee-41d0-9acc-ae7230400019"
    "D66" = "53af4926-52ee-41d0-9acc-ae7230400020"
    "D67" = "53af4926-52ee-41d0-9acc-ae7230400016"
    "D69" = "53af4926-52ee-41d0-9acc-ae7230400018"
    "D70" = "53af4926-52ee-41d0-9acc-ae7230400019"
    "D71" = "53af4926-52ee-41d0-9acc-ae7230400020"
    "D72" = "53af4926-52ee-41d0-9acc-ae7230400016"
    "D73" = "53af4926-52ee-41d0-9acc-ae7230400017"
    "D74" = "53af4926-52ee-41d0-9acc-ae7230400018"
    "D76" = "53af4926-52ee-41d0-9acc-ae7230400020"
    "D77" = "53af4926-52ee-41d0-9acc-ae7230400016"
    "D78" = "53af4926-52ee-41d0-9acc-ae7230400017"
    "D79" = "53af4926-52ee-41d0-9acc-ae7230400018"
    "D80" = "53af4926-52ee-41d0-9acc-ae7230400019"
    "D83" = "53af4926-52ee-41d0-9acc-ae7230400022"
    "D84" = "53af4926-52ee-41d0-9acc-ae7230400023"
    "D85" = "53af4926-52ee-41d0-9acc-ae7230400024"
    "D86" = "53af4926-52ee-41d0-9acc-ae7230400025"
    "D87" = "53af4926-52ee-41d0-9acc-ae7230400021"
    "D89" = "53af4926-52ee-41d0-9acc-ae7230400023"
    "D90" = "53af4926-52ee-41d0-9acc-ae7230400024"
    "D91" = "53af4926-52ee-41d0-9acc-ae7230400025"
    "D92" = "53af4926-52ee-41d0-9acc-ae7230400021"
    "D93" = "53af4926-52ee-41d0-9acc-ae7230400022"
    "D94" = "53af4926-52ee-41d0-9acc-ae7230400023"
    "D96" = "53af4926-52ee-41d0-9acc-ae7230400025"
    "D97" = "53af4926-52ee-41d0-9acc-ae7230400021"
    "D98" = "53af4926-52ee-41d0-9acc-ae7230400022"
    "D99" = "53af4926-52ee-41d0-9acc-ae7230400023"
    "D100" = "53af4926-52ee-41d0-9acc-ae7230400024"
    "D103" = "53af4926-52ee-41d0-9acc-ae7230400027"
    "D104" = "53af4926-52ee-41d0-9acc-ae7230400028"
    "D105" = "53af4926-52ee-41d0-9acc-ae7230400029"
    "D106" = "53af4926-52ee-41d0-9acc-ae7230400030"
    "D107" = "53af4926-52ee-41d0-9acc-ae7230400026"
    "D109" = "53af4926-52ee-41d0-9acc-ae7230400028"
    "D110" = "53af4926-52ee-41d0-9acc-ae7230400029"
    "D111" = "53af4926-52ee-41d0-9acc-ae7230400030"
    "D112" = "53af4926-52ee-41d0-9acc-ae7230400026"
    "D113" = "53af4926-52ee-41d0-9acc-ae7230400027"
    "D114" = "53af4926-52ee-41d0-9acc-ae7230400028"
    "D116" = "53af4926-52ee-41d0-9acc-ae7230400030"
    "D117" = "53af4926-52ee-41d0-9acc-ae7230400026"
    "D118" = "53af4926-52ee-41d0-9acc-ae7230400027"
    "D119" = "53af4926-52ee-41d0-9acc-ae7230400028"
    "D120" = "53af4926-52ee-41d0-9acc-ae7230400029"
    "D123" = "53af4926-52ee-41d0-9acc-ae7230400032"
    "D124" = "53af4926-52ee-41d0-9acc-ae7230400033"
    "D125" = "53af4926-52ee-41d0-9acc-ae7230400034"
    "D126" = "53af4926-52ee-41d0-9acc-ae7230400035"
    "D127" = "53af4926-52ee-41d0-9acc-ae7230400031"
    "D129" = "53af4926-52ee-41d0-9acc-ae7230400033"
    "D130" = "53af4926-52ee-41d0-9acc-ae7230400034"
    "D131" = "53af4926-52ee-41d0-9acc-ae7230400035"
    "D132" = "53af4926-52ee-41d0-9acc-ae7230400031"
    "D133" = "53af4926-52ee-41d0-9acc-ae7230400032"
    "D134" = "53af4926-52ee-41d0-9acc-ae7230400033"
    "D136" = "53af4926-52ee-41d0-9acc-ae7230400035"
    "D137" = "53af4926-52ee-41d0-9acc-ae7230400031"
    "D138" = "53af4926-52ee-41d0-9acc-ae7230400032"
    "D139" = "53af4926-52ee-41d0-9acc-ae7230400033"
    "D140" = "53af4926-52ee-41d0-9acc-ae7230400034"
    "D143" = "53af4926-52ee-41d0-9acc-ae7230400037"
    "D144" = "53af4926-52ee-41d0-9acc-ae7230400038"
    "D145" = "53af4926-52ee-41d0-9acc-ae7230400039"
    "D146" = "53af4926-52ee-41d0-9acc-ae7230400040"
    "D147" = "53af4926-52ee-41d0-9acc-ae7230400036"
    "D149" = "53af4926-52ee-41d0-9acc-ae7230400038"
    "D150" = "53af4926-52ee-41d0-9acc-ae7230400039"
    "D151" = "53af4926-52ee-41d0-9acc-ae7230400040"
    "D152" = "53af4926-52ee-41d0-9acc-ae7230400036"
    "D153" = "53af4926-52ee-41d0-9acc-ae7230400037"
    "D154" = "53af4926-52ee-41d0-9acc-ae7230400038"
    "D156" = "53af4926-52ee-41d0-9acc-ae7230400040"
    "D157" = "53af4926-52ee-41d0-9acc-ae7230400036"
    "D158" = "53af4926-52ee-41d0-9acc-ae7230400037"
    "D159" = "53af4926-52ee-41d0-9acc-ae7230400038"
    "D160" = "53af4926-52ee-41d0-9acc-ae7230400039"
    "D163" = "53af4926-52ee-41d0-9acc-ae7230400042"
    "D164" = "53af4926-52ee-41d0-9acc-ae7230400043"
    "D165" = "53af4926-52ee-41d0-9acc-ae7230400044"
    "D166" = "53af4926-52ee-41d0-9acc-ae7230400045"
    "D167" = "53af4926-52ee-41d0-9acc-ae7230400041"
    "D169" = "53af4926-52ee-41d0-9acc-ae7230400043"
    "D170" = "53af4926-52ee-41d0-9acc-ae7230400044"
    "D171" = "53af4926-52ee-41d0-9acc-ae7230400045"
    "D172" = "53af4926-52ee-41d0-9acc-ae7230400041"
    "D173" = "53af4926-52ee-41d0-9acc-ae7230400042"
    "D174" = "53af4926-52ee-41d0-9acc-ae7230400043"
    "D176" = "53af4926-52ee-41d0-9acc-ae7230400045"
    "D177" = "53af4926-52ee-41d0-9acc-ae7230400041"
    "D178" = "53af4926-52ee-41d0-9acc-ae7230400042"
    "D179" = "53af4926-52ee-41d0-9acc-ae7230400043"
    "D180" = "53af4926-52ee-41d0-9acc-ae7230400044"
    "D183" = "53af4926-52ee-41d0-9acc-ae7230400047"
    "D184" = "53af4926-52ee-41d0-9acc-ae7230400048"
    "D185" = "53af4926-52ee-41d0-9acc-ae7230400049"
    "D186" = "53af4926-52ee-41d0-9acc-ae7230400050"
    "D187" = "53af4926-52ee-41d0-9acc-ae7230400046"
    "D189" = "53af4926-52ee-41d0-9acc-ae7230400048"
    "D190" = "53af4926-52ee-41d0-9acc-ae7230400049"
    "D191" = "53af4926-52ee-41d0-9acc-ae7230400050"
    "D192" = "53af4926-52ee-41d0-9acc-ae7230400046"
    "D193" = "53af4926-52ee-41d0-9acc-ae7230400047"
    "D194" = "53af4926-52ee-41d0-9acc-ae7230400048"
    "D196" = "53af4926-52ee-41d0-9acc-ae7230400050"
    "D197" = "53af4926-52ee-41d0-9acc-ae7230400046"
    "D198" = "53af4926-52ee-41d0-9acc-ae7230400047"
    "D199" = "53af4926-52ee-41d0-9acc-ae7230400048"
    "D200" = "53af4926-52ee-41d0-9acc-ae7230400049"
    "D203" = "53af4926-52ee-41d0-9acc-ae7230400052"
    "D204" = "53af4926-52ee-41d0-9acc-ae7230400053"
    "D205" = "53af4926-52ee-41d0-9acc-ae7230400054"
    "D206" = "53af4926-52ee-41d0-9acc-ae7230400055"
    "D207" = "53af4926-52ee-41d0-9acc-ae7230400051"
    "D209" = "53af4926-52ee-41d0-9acc-ae7230400053"
    "D210" = "53af4926-52ee-41d0-9acc-ae7230400054"
    "D211" = "53af4926-52ee-41d0-9acc-ae7230400055"
    "D212" = "53af4926-52ee-41d0-9acc-ae7230400051"
    "D213" = "53af4926-52ee-41d0-9acc-ae7230400052"
    "D214" = "53af4926-52ee-41d0-9acc-ae7230400053"
    "D216" = "53af4926-52ee-41d0-9acc-ae7230400055"
    "D217" = "53af4926-52ee-41d0-9acc-ae7230400051"
    "D218" = "53af4926-52ee-41d0-9acc-ae7230400052"
    "D219" = "53af4926-52ee-41d0-9acc-ae7230400053"
    "D220" = "53af4926-52ee-41d0-9acc-ae7230400054"
    "D223" = "53af4926-52ee-41d0-9acc-ae7230400057"
    "D224" = "53af4926-52ee-41d0-9acc-ae7230400058"
    "D225" = "53af4926-52ee-41d0-9acc-ae7230400059"
    "D226" = "53af4926-52ee-41d0-9acc-ae7230400060"
    "D227" = "53af4926-52ee-41d0-9acc-ae7230400056"
    "D229" = "53af4926-52ee-41d0-9acc-ae7230400058"
    "D230" = "53af4926-52ee-41d0-9acc-ae7230400059"
    "D231" = "53af4926-52ee-41d0-9acc-ae7230400060"
    "D232" = "53af4926-52ee-41d0-9acc-ae7230400056"
    "D233" = "53af4926-52ee-41d0-9acc-ae7230400057"
    "D234" = "53af4926-52ee-41d0-9acc-ae7230400058"
    "D236" = "53af4926-52ee-41d0-9acc-ae7230400060"
    "D237" = "53af4926-52ee-41d0-9acc-ae7230400056"
    "D238" = "53af4926-52ee-41d0-9acc-ae7230400057"
    "D239" = "53af4926-52ee-41d0-9acc-ae7230400058"
    "D240" = "53af4926-52ee-41d0-9acc-ae7230400059"
    "D243" = "53af4926-52ee-41d0-9acc-ae7230400062"
    "D244" = "53af4926-52ee-41d0-9acc-ae7230400063"
    "D245" = "53af4926-52ee-41d0-9acc-ae7230400064"
    "D246" = "53af4926-52ee-41d0-9acc-ae7230400065"
    "D247" = "53af4926-52ee-41d0-9acc-ae7230400061"
    "D249" = "53af4926-52ee-41d0-9acc-ae7230400063"
    "D250" = "53af4926-52ee-41d0-9acc-ae7230400064"
    "D251" = "53af4926-52ee-41d0-9acc-ae7230400065"
    "D252" = "53af4926-52ee-41d0-9acc-ae7230400061"
    "D253" = "53af4926-52ee-41d0-9acc-ae7230400062"
    "D254" = "53af4926-52ee-41d0-9acc-ae7230400063"
    "D256" = "53af4926-52ee-41d0-9acc-ae7230400065"
    "D257" = "53af4926-52ee-41d0-9acc-ae7230400061"
    "D258" = "53af4926-52ee-41d0-9acc-ae7230400062"
    "D259" = "53af4926-52ee-41d0-9acc-ae7230400063"
    "D260" = "53af4926-52ee-41d0-9acc-ae7230400064"
    "D263" = "53af4926-52ee-41d0-9acc-ae7230400067"
    "D264" = "53af4926-52ee-41d0-9acc-ae7230400068"
    "D265" = "53af4926-52ee-41d0-9acc-ae7230400069"
    "D266" = "53af4926-52ee-41d0-9acc-ae7230400070"
    "D267" = "53af4926-52ee-41d0-9acc-ae7230400066"
    "D269" = "53af4926-52ee-41d0-9acc-ae7230400068"
    "D270" = "53af4926-52ee-41d0-9acc-ae7230400069"
    "D271" = "53af4926-52ee-41d0-9acc-ae7230400070"
    "D272" = "53af4926-52ee-41d0-9acc-ae7230400066"
    "D273" = "53af4926-52ee-41d0-9acc-ae7230400067"
    "D274" = "53af4926-52ee-41d0-9acc-ae7230400068"
    "D276" = "53af4926-52ee-41d0-9acc-ae7230400070"
    "D277" = "53af4926-52ee-41d0-9acc-ae7230400066"
    "D278" = "53af4926-52ee-41d0-9acc-ae7230400067"
    "D279" = "53af4926-52ee-41d0-9acc-ae7230400068"
    "D280" = "53af4926-52ee-41d0-9acc-ae7230400069"
    "D283" = "53af4926-52ee-41d0-9acc-ae7230400072"
    "D284" = "53af4926-52ee-41d0-9acc-ae7230400073"
    "D285" = "53af4926-52ee-41d0-9acc-ae7230400074"
    "D286" = "53af4926-52ee-41d0-9acc-ae7230400075"
    "D287" = "53af4926-52ee-41d0-9acc-ae7230400071"
    "D289" = "53af4926-52ee-41d0-9acc-ae7230400073"
    "D290" = "53af4926-52ee-41d0-9acc-ae7230400074"
    "D291" = "53af4926-52ee-41d0-9acc-ae7230400075"
    "D292" = "53af4926-52ee-41d0-9acc-ae7230400071"
    "D293" = "53af4926-52ee-41d0-9acc-ae7230400072"
    "D294" = "53af4926-52ee-41d0-9acc-ae7230400073"
    "D296" = "53af4926-52ee-41d0-9acc-ae7230400075"
    "D297" = "53af4926-52ee-41d0-9acc-ae7230400071"
    "D298" = "53af4926-52ee-41d0-9acc-ae7230400072"
    "D299" = "53af4926-52ee-41d0-9acc-ae7230400073"
    "D300" = "53af4926-52ee-41d0-9acc-ae7230400074"
    "D303" = "53af4926-52ee-41d0-9acc-ae7230400077"
    "D304" = "53af4926-52ee-41d0-9acc-ae7230400078"
    "D305" = "53af4926-52ee-41d0-9acc-ae7230400079"
    "D306" = "53af4926-52ee-41d0-9acc-ae7230400080"
    "D307" = "53af4926-52ee-41d0-9acc-ae7230400076"
    "D309" = "53af4926-52ee-41d0-9acc-ae7230400078"
    "D310" = "53af4926-52ee-41d0-9acc-ae7230400079"
    "D311" = "53af4926-52ee-41d0-9acc-ae7230400080"
    "D312" = "53af4926-52ee-41d0-9acc-ae7230400076"
    "D313" = "53af4926-52ee-41d0-9acc-ae7230400077"
    "D314" = "53af4926-52ee-41d0-9acc-ae7230400078"
    "D316" = "53af4926-52ee-41d0-9acc-ae7230400080"
    "D317" = "53af4926-52ee-41d0-9acc-ae7230400076"
    "D318" = "53af4926-52ee-41d0-9acc-ae7230400077"
    "D319" = "53af4926-52ee-41d0-9acc-ae7230400078"
    "D320" = "53af4926-52ee-41d0-9acc-ae7230400079"
    "D323" = "53af4926-52ee-41d0-9acc-ae7230400082"
    "D324" = "53af4926-52ee-41d0-9acc-ae7230400083"
    "D325" = "53af4926-52ee-41d0-9acc-ae7230400084"
    "D326" = "53af4926-52ee-41d0-9acc-ae7230400085"
    "D327" = "53af4926-52ee-41d0-9acc-ae7230400081"
    "D329" = "53af4926-52ee-41d0-9acc-ae7230400083"
    "D330" = "53af4926-52ee-41d0-9acc-ae7230400084"
    "D331" = "53af4926-52ee-41d0-9acc-ae7230400085"
    "D332" = "53af4926-52ee-41d0-9acc-ae7230400081"
    "D333" = "53af4926-52ee-41d0-9acc-ae7230400082"
    "D334" = "53af4926-52ee-41d0-9acc-ae7230400083"
    "D336" = "53af4926-52ee-41d0-9acc-ae7230400085"
    "D337" = "53af4926-52ee-41d0-9acc-ae7230400081"
    "D338" = "53af4926-52ee-41d0-9acc-ae7230400082"
    "D339" = "53af4926-52ee-41d0-9acc-ae7230400083"
    "D340" = "53af4926-52ee-41d0-9acc-ae7230400084"
    "D343" = "53af4926-52ee-41d0-9acc-ae7230400087"
    "D344" = "53af4926-52ee-41d0-9acc-ae7230400088"
    "D345" = "53af4926-52ee-41d0-9acc-ae7230400089"
    "D346" = "53af4926-52ee-41d0-9acc-ae7230400090"
    "D347" = "53af4926-52ee-41d0-9acc-ae7230400086"
    "D349" = "53af4926-52ee-41d0-9acc-ae7230400088"
    "D350" = "53af4926-52ee-41d0-9acc-ae7230400089"
    "D351" = "53af4926-52ee-41d0-9acc-ae7230400090"
    "D352" = "53af4926-52ee-41d0-9acc-ae7230400086"
    "D353" = "53af4926-52ee-41d0-9acc-ae7230400087"
    "D354" = "53af4926-52ee-41d0-9acc-ae7230400088"
    "D356" = "53af4926-52ee-41d0-9acc-ae7230400090"
    "D357" = "53af4926-52ee-41d0-9acc-ae7230400086"
    "D358" = "53af4926-52ee-41d0-9acc-ae7230400087"
    "D359" = "53af4926-52ee-41d0-9acc-ae7230400088"
    "D360" = "53af4926-52ee-41d0-9acc-ae7230400089"
    "D363" = "53af4926-52ee-41d0-9acc-ae7230400092"
    "D364" = "53af4926-52ee-41d0-9acc-ae7230400093"
    "D365" = "53af4926-52ee-41d0-9acc-ae7230400094"
    "D366" = "53af4926-52ee-41d0-9acc-ae7230400095"
    "D367" = "53af4926-52ee-41d0-9acc-ae7230400091"
    "D369" = "53af4926-52ee-41d0-9acc-ae7230400093"
    "D370" = "53af4926-52ee-41d0-9acc-ae7230400094"
    "D371" = "53af4926-52ee-41d0-9acc-ae7230400095"
    "D372" = "53af4926-52ee-41d0-9acc-ae7230400091"
    "D373" = "53af4926-52ee-41d0-9acc-ae7230400092"
    "D374" = "53af4926-52ee-41d0-9acc-ae7230400093"
    "D376" = "53af4926-52ee-41d0-9acc-ae7230400095"
    "D377" = "53af4926-52ee-41d0-9acc-ae7230400091"
    "D378" = "53af4926-52ee-41d0-9acc-ae7230400092"
    "D379" = "53af4926-52ee-41d0-9acc-ae7230400093"
    "D380" = "53af4926-52ee-41d0-9acc-ae7230400094"
    "D383" = "53af4926-52ee-41d0-9acc-ae7230400097"
    "D384" = "53af4926-52ee-41d0-9acc-ae7230400098"
    "D385" = "53af4926-52ee-41d0-9acc-ae7230400099"
    "D386" = "53af4926-52ee-41d0-9acc-ae7230400100"
    "D387" = "53af4926-52ee-41d0-9acc-ae7230400096"
    "D389" = "53af4926-52ee-41d0-9acc-ae7230400098"
    "D390" = "53af4926-52ee-41d0-9acc-ae7230400099"
    "D391" = "53af4926-52ee-41d0-9acc-ae7230400100"
    "D392" = "53af4926-52ee-41d0-9acc-ae7230400096"
    "D393" = "53af4926-52ee-41d0-9acc-ae7230400097"
    "D394" = "53af4926-52ee-41d0-9acc-ae7230400098"
    "D396" = "53af4926-52ee-41d0-9acc-ae7230400100"
    "D397" = "53af4926-52ee-41d0-9acc-ae7230400096"
    "D398" = "53af4926-52ee-41d0-9acc-ae7230400097"
    "D399" = "53af4926-52ee-41d0-9acc-ae7230400098"
    "D400" = "53af4926-52ee-41d0-9acc-ae7230400099"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

[void]$ws.Range("A14").Select()

Write-Host "Updated $($updates.Count) teacher_id cells"